$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C1").Value = 10

$ws.Range("C3").Value = "Carlo Zine Marc Sam"
$ws.Range("C4").Value = "Rief , Michiel , Robin"

$ws.Range("C8").Value = "Player Controller"
$ws.Range("C9").Value = "unwrap pistol"
$ws.Range("C10").Value = "main menu"
$ws.Range("C11").Value = "Main Character"
$ws.Range("C12").Value = "Muren en Textures"

$ws.Range("C14").Value = "Shotgun texture"
$ws.Range("C15").Value = "Planning volgen"
$ws.Range("C16").Value = "Planning Maken"
$ws.Range("C17").Value = "Texture en Animatie"
$ws.Range("C18").Value = "Planning Maken"
$ws.Range("C19").Value = "Main Character]"
$ws.Range("C20").Value = "Muren en Textures"

$ws.Range("C22").Value = ""
$ws.Range("C23").Value = "Waarom zo boos?(-zine)"
$ws.Range("C24").Value = "Anxiety"

$ws.Range("C29").Value = "Team Spirit"
$ws.Range("C31").Value = ""
